$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOGINDATA")

# Add new row 5 data first, so "addCoupon" gets added to the shared
# string table before "no" does.
$ws.Range("A5").Value = "addCoupon"
$ws.Range("B5").Value = "yes"
$ws.Range("C5").Value = "chrome"
$ws.Range("D5").Value = "'"
$ws.Range("E5").Value = "'"

# Update existing rows: change "yes" -> "no" in column B for rows 2-4
$ws.Range("B2").Value = "no"
$ws.Range("B3").Value = "no"
$ws.Range("B4").Value = "no"

# Update selection to B4 to match the new active cell
$ws.Range("B4").Select()
